$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 582.0909 # H2
$ws.Cells.Item(2, 9).Value = 289.22223 # I2
$ws.Cells.Item(2, 11).Value = 289.22223 # K2
$ws.Cells.Item(2, 13).Value = -176.22223 # M2
$ws.Cells.Item(10, 8).Value = 7926.3335 # H10
$ws.Cells.Item(10, 10).Value = 7926.3335 # J10
$ws.Cells.Item(10, 12).Value = 7926.3335 # L10
$ws.Cells.Item(10, 14).Value = -8512.333500000001 # N10
$ws.Cells.Item(134, 8).Value = 172040.5 # H134
$ws.Cells.Item(134, 10).Value = 192448.6 # J134
$ws.Cells.Item(134, 12).Value = 192448.6 # L134
$ws.Cells.Item(134, 14).Value = -202588.6 # N134
$ws.Cells.Item(138, 8).Value = 3188.054 # H138
$ws.Cells.Item(138, 9).Value = 2013.6666 # I138
$ws.Cells.Item(138, 10).Value = 3988.7727 # J138
$ws.Cells.Item(138, 11).Value = 6040.9998 # K138
$ws.Cells.Item(138, 12).Value = 11966.3181 # L138
$ws.Cells.Item(138, 13).Value = -900.9997999999996 # M138
$ws.Cells.Item(138, 14).Value = -22246.3181 # N138

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2020762.5 # H2
$ws.Cells.Item(2, 9).Value = 2828232 # I2
$ws.Cells.Item(2, 10).Value = 2089 # J2
$ws.Cells.Item(2, 11).Value = 2828232 # K2
$ws.Cells.Item(2, 12).Value = 2089 # L2
$ws.Cells.Item(2, 13).Value = -2828119 # M2
$ws.Cells.Item(2, 14).Value = -2315 # N2
$ws.Cells.Item(45, 8).Value = 6692344.5 # H45
$ws.Cells.Item(45, 9).Value = 9617327 # I45
$ws.Cells.Item(45, 11).Value = 9617327 # K45
$ws.Cells.Item(45, 13).Value = -9616950 # M45
$ws.Cells.Item(61, 8).Value = 5681.816 # H61
$ws.Cells.Item(61, 9).Value = 6060.2812 # I61
$ws.Cells.Item(61, 10).Value = 3663.3333 # J61
$ws.Cells.Item(61, 11).Value = 6060.2812 # K61
$ws.Cells.Item(61, 12).Value = 3663.3333 # L61
$ws.Cells.Item(61, 13).Value = -5848.2812 # M61
$ws.Cells.Item(61, 14).Value = -4087.3333 # N61
$ws.Cells.Item(74, 8).Value = 99684.95 # H74
$ws.Cells.Item(74, 9).Value = 87744.836 # I74
$ws.Cells.Item(74, 10).Value = 117595.125 # J74
$ws.Cells.Item(74, 11).Value = 87744.836 # K74
$ws.Cells.Item(74, 12).Value = 117595.125 # L74
$ws.Cells.Item(74, 13).Value = -86870.836 # M74
$ws.Cells.Item(74, 14).Value = -119343.125 # N74
$ws.Cells.Item(77, 8).Value = 99684.95 # H77
$ws.Cells.Item(77, 9).Value = 87744.836 # I77
$ws.Cells.Item(77, 10).Value = 117595.125 # J77
$ws.Cells.Item(77, 11).Value = 438724.18 # K77
$ws.Cells.Item(77, 12).Value = 587975.625 # L77
$ws.Cells.Item(77, 13).Value = -434356.18 # M77
$ws.Cells.Item(77, 14).Value = -596711.625 # N77
$ws.Cells.Item(116, 8).Value = 2020762.5 # H116
$ws.Cells.Item(116, 9).Value = 2828232 # I116
$ws.Cells.Item(116, 10).Value = 2089 # J116
$ws.Cells.Item(116, 11).Value = 2828232 # K116
$ws.Cells.Item(116, 12).Value = 2089 # L116
$ws.Cells.Item(116, 13).Value = -2825938 # M116
$ws.Cells.Item(116, 14).Value = -6677 # N116
$ws.Cells.Item(136, 8).Value = 5681.816 # H136
$ws.Cells.Item(136, 9).Value = 6060.2812 # I136
$ws.Cells.Item(136, 10).Value = 3663.3333 # J136
$ws.Cells.Item(136, 11).Value = 18180.8436 # K136
$ws.Cells.Item(136, 12).Value = 10989.9999 # L136
$ws.Cells.Item(136, 13).Value = -15630.8436 # M136
$ws.Cells.Item(136, 14).Value = -16089.9999 # N136
$ws.Cells.Item(139, 8).Value = 59999.5 # H139
$ws.Cells.Item(139, 10).Value = 59999.5 # J139
$ws.Cells.Item(139, 12).Value = 59999.5 # L139
$ws.Cells.Item(139, 14).Value = -70279.5 # N139

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2020762.5 # H3
$ws.Cells.Item(3, 9).Value = 2828232 # I3
$ws.Cells.Item(3, 10).Value = 2089 # J3
$ws.Cells.Item(3, 11).Value = 2828232 # K3
$ws.Cells.Item(3, 12).Value = 2089 # L3
$ws.Cells.Item(3, 13).Value = -2828118 # M3
$ws.Cells.Item(3, 14).Value = -2317 # N3
$ws.Cells.Item(20, 8).Value = 22225052 # H20
$ws.Cells.Item(20, 9).Value = 25642444 # I20
$ws.Cells.Item(20, 11).Value = 25642444 # K20
$ws.Cells.Item(20, 13).Value = -25642197 # M20
$ws.Cells.Item(86, 8).Value = 5564679.5 # H86
$ws.Cells.Item(86, 9).Value = 11114409 # I86
$ws.Cells.Item(86, 11).Value = 11114409 # K86
$ws.Cells.Item(86, 13).Value = -11113286 # M86
$ws.Cells.Item(89, 8).Value = 5564679.5 # H89
$ws.Cells.Item(89, 9).Value = 11114409 # I89
$ws.Cells.Item(89, 11).Value = 55572045 # K89
$ws.Cells.Item(89, 13).Value = -55566429 # M89
$ws.Cells.Item(94, 8).Value = 2781561.2 # H94
$ws.Cells.Item(94, 9).Value = 3031476 # I94
$ws.Cells.Item(94, 11).Value = 3031476 # K94
$ws.Cells.Item(94, 13).Value = -3031025 # M94
$ws.Cells.Item(105, 8).Value = 5683135.5 # H105
$ws.Cells.Item(105, 9).Value = 5683135.5 # I105
$ws.Cells.Item(105, 10).Value = 0 # J105
$ws.Cells.Item(105, 11).Value = 5683135.5 # K105
$ws.Cells.Item(105, 12).Value = 0 # L105
$ws.Cells.Item(105, 13).Value = -5681388.5 # M105
$ws.Cells.Item(105, 14).ClearContents() # N105
$ws.Cells.Item(134, 8).Value = 11681.767 # H134
$ws.Cells.Item(134, 9).Value = 9194.120000000001 # I134
$ws.Cells.Item(134, 11).Value = 27582.36 # K134
$ws.Cells.Item(134, 13).Value = -25047.36 # M134

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(11, 8).Value = 5666 # H11
$ws.Cells.Item(11, 9).Value = 0 # I11
$ws.Cells.Item(11, 10).Value = 5666 # J11
$ws.Cells.Item(11, 11).Value = 0 # K11
$ws.Cells.Item(11, 12).Value = 5666 # L11
$ws.Cells.Item(11, 13).ClearContents() # M11
$ws.Cells.Item(11, 14).Value = -5946 # N11
$ws.Cells.Item(31, 8).Value = 26837.861 # H31
$ws.Cells.Item(31, 9).Value = 6797.5557 # I31
$ws.Cells.Item(31, 11).Value = 6797.5557 # K31
$ws.Cells.Item(31, 13).Value = -6502.5557 # M31
$ws.Cells.Item(34, 8).Value = 26837.861 # H34
$ws.Cells.Item(34, 9).Value = 6797.5557 # I34
$ws.Cells.Item(34, 11).Value = 6797.5557 # K34
$ws.Cells.Item(34, 13).Value = -6595.5557 # M34
$ws.Cells.Item(97, 8).Value = 0 # H97
$ws.Cells.Item(97, 10).Value = 0 # J97
$ws.Cells.Item(97, 12).Value = 0 # L97
$ws.Cells.Item(97, 14).ClearContents() # N97
$ws.Cells.Item(109, 8).Value = 23164.572 # H109
$ws.Cells.Item(109, 10).Value = 23164.572 # J109
$ws.Cells.Item(109, 12).Value = 23164.572 # L109
$ws.Cells.Item(109, 14).Value = -25244.572 # N109
$ws.Cells.Item(134, 8).Value = 41533.793 # H134
$ws.Cells.Item(134, 9).Value = 52634.25 # I134
$ws.Cells.Item(134, 10).Value = 16866.111 # J134
$ws.Cells.Item(134, 11).Value = 157902.75 # K134
$ws.Cells.Item(134, 12).Value = 50598.333 # L134
$ws.Cells.Item(134, 13).Value = -155367.75 # M134
$ws.Cells.Item(134, 14).Value = -55668.333 # N134

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(18, 8).Value = 467.8889 # H18
$ws.Cells.Item(18, 9).Value = 261.83334 # I18
$ws.Cells.Item(18, 11).Value = 785.5000200000001 # K18
$ws.Cells.Item(18, 13).Value = -616.5000200000001 # M18
$ws.Cells.Item(58, 8).Value = 1692.1538 # H58
$ws.Cells.Item(58, 9).Value = 999 # I58
$ws.Cells.Item(58, 11).Value = 2997 # K58
$ws.Cells.Item(58, 13).Value = -2869 # M58
$ws.Cells.Item(107, 8).Value = 1415.7693 # H107
$ws.Cells.Item(107, 9).Value = 2296 # I107
$ws.Cells.Item(107, 11).Value = 6888 # K107
$ws.Cells.Item(107, 13).Value = -4968 # M107
$ws.Cells.Item(131, 8).Value = 16671388 # H131
$ws.Cells.Item(131, 9).Value = 5954019.5 # I131
$ws.Cells.Item(131, 10).Value = 30311676 # J131
$ws.Cells.Item(131, 11).Value = 17862058.5 # K131
$ws.Cells.Item(131, 12).Value = 90935028 # L131
$ws.Cells.Item(131, 13).Value = -17857018.5 # M131
$ws.Cells.Item(131, 14).Value = -90945108 # N131

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 200000000 # H70
$ws.Cells.Item(70, 9).Value = 200000000 # I70
$ws.Cells.Item(70, 10).Value = 0 # J70
$ws.Cells.Item(70, 11).Value = 200000000 # K70
$ws.Cells.Item(70, 12).Value = 0 # L70
$ws.Cells.Item(70, 13).Value = -199999730 # M70
$ws.Cells.Item(70, 14).ClearContents() # N70
$ws.Cells.Item(73, 8).Value = 200000000 # H73
$ws.Cells.Item(73, 9).Value = 200000000 # I73
$ws.Cells.Item(73, 10).Value = 0 # J73
$ws.Cells.Item(73, 11).Value = 200000000 # K73
$ws.Cells.Item(73, 12).Value = 0 # L73
$ws.Cells.Item(73, 13).Value = -199999064 # M73
$ws.Cells.Item(73, 14).ClearContents() # N73

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 5008749.5 # H2
$ws.Cells.Item(2, 9).Value = 10000000 # I2
$ws.Cells.Item(2, 11).Value = 10000000 # K2
$ws.Cells.Item(2, 13).Value = -9999888 # M2
$ws.Cells.Item(100, 8).Value = 3550 # H100
$ws.Cells.Item(100, 9).Value = 2000 # I100
$ws.Cells.Item(100, 10).Value = 5381.8184 # J100
$ws.Cells.Item(100, 11).Value = 2000 # K100
$ws.Cells.Item(100, 12).Value = 5381.8184 # L100
$ws.Cells.Item(100, 13).Value = -1459 # M100
$ws.Cells.Item(100, 14).Value = -6463.8184 # N100
$ws.Cells.Item(125, 8).Value = 71497.5 # H125
$ws.Cells.Item(125, 10).Value = 71497.5 # J125
$ws.Cells.Item(125, 12).Value = 71497.5 # L125
$ws.Cells.Item(125, 14).Value = -81337.5 # N125

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 47620436 # H107
$ws.Cells.Item(107, 9).Value = 55556120 # I107
$ws.Cells.Item(107, 11).Value = 166668360 # K107
$ws.Cells.Item(107, 13).Value = -166666440 # M107
$ws.Cells.Item(126, 8).Value = 2425.7144 # H126
$ws.Cells.Item(126, 9).Value = 2513.75 # I126
$ws.Cells.Item(126, 11).Value = 7541.25 # K126
$ws.Cells.Item(126, 13).Value = -5071.25 # M126
$ws.Cells.Item(136, 8).Value = 8190.5415 # H136
$ws.Cells.Item(136, 9).Value = 13253.385 # I136
$ws.Cells.Item(136, 10).Value = 2207.182 # J136
$ws.Cells.Item(136, 11).Value = 39760.155 # K136
$ws.Cells.Item(136, 12).Value = 6621.545999999999 # L136
$ws.Cells.Item(136, 13).Value = -37210.155 # M136
$ws.Cells.Item(136, 14).Value = -11721.546 # N136
